# Apply scheduled-runner price/profit updates to Sheets/Excalibur_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 771.1429000000001
$ws.Range("I33").Value = 682.6667
$ws.Range("J33").Value = 837.5
$ws.Range("K33").Value = 682.6667
$ws.Range("L33").Value = 837.5
$ws.Range("M33").Value = -453.6667
$ws.Range("N33").Value = -1295.5

$ws.Range("H92").Value = 4396.2
$ws.Range("I92").Value = 5244.625
$ws.Range("K92").Value = 5244.625
$ws.Range("M92").Value = -3996.625

$ws.Range("H98").Value = 1345.4595
$ws.Range("I98").Value = 1345.4595
$ws.Range("K98").Value = 1345.4595
$ws.Range("M98").Value = 152.5405000000001

$ws.Range("H112").Value = 1435.7037
$ws.Range("I112").Value = 874.3333
$ws.Range("J112").Value = 1505.875
$ws.Range("K112").Value = 2622.9999
$ws.Range("L112").Value = 4517.625
$ws.Range("M112").Value = -1514.9999
$ws.Range("N112").Value = -6733.625

$ws.Range("H122").Value = 1345.4595
$ws.Range("I122").Value = 1345.4595
$ws.Range("K122").Value = 4036.3785
$ws.Range("M122").Value = -1586.3785

$ws.Range("H132").Value = 1560.6604
$ws.Range("I132").Value = 1312.4694
$ws.Range("J132").Value = 4601
$ws.Range("K132").Value = 3937.4082
$ws.Range("L132").Value = 13803
$ws.Range("M132").Value = -1407.4082
$ws.Range("N132").Value = -18863

$ws.Range("H138").Value = 3205.4
$ws.Range("J138").Value = 3983.9285
$ws.Range("L138").Value = 11951.7855
$ws.Range("N138").Value = -22231.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3557.38
$ws.Range("I32").Value = 2566.5803
$ws.Range("K32").Value = 2566.5803
$ws.Range("M32").Value = -2279.5803

$ws.Range("H45").Value = 4956.724
$ws.Range("I45").Value = 4342.3335
$ws.Range("J45").Value = 6569.5
$ws.Range("K45").Value = 4342.3335
$ws.Range("L45").Value = 6569.5
$ws.Range("M45").Value = -3965.3335
$ws.Range("N45").Value = -7323.5

$ws.Range("H61").Value = 10376.704
$ws.Range("I61").Value = 10617
$ws.Range("K61").Value = 10617
$ws.Range("M61").Value = -10405

$ws.Range("H122").Value = 2311.3667
$ws.Range("I122").Value = 1130.1666
$ws.Range("K122").Value = 3390.4998
$ws.Range("M122").Value = -940.4998000000001

$ws.Range("H136").Value = 10376.704
$ws.Range("I136").Value = 10617
$ws.Range("K136").Value = 31851
$ws.Range("M136").Value = -29301

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26371.773
$ws.Range("I134").Value = 29740.805
$ws.Range("K134").Value = 89222.41500000001
$ws.Range("M134").Value = -86687.41500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 890.86664
$ws.Range("I22").Value = 665.9091
$ws.Range("J22").Value = 1509.5
$ws.Range("K22").Value = 665.9091
$ws.Range("L22").Value = 1509.5
$ws.Range("M22").Value = -315.9091
$ws.Range("N22").Value = -2209.5

$ws.Range("H105").Value = 876.5714
$ws.Range("I105").Value = 876.5714
$ws.Range("K105").Value = 876.5714
$ws.Range("M105").Value = 870.4286

$ws.Range("H122").Value = 1962.875
$ws.Range("I122").Value = 1430.3
$ws.Range("J122").Value = 2850.5
$ws.Range("K122").Value = 4290.9
$ws.Range("L122").Value = 8551.5
$ws.Range("M122").Value = -1840.9
$ws.Range("N122").Value = -13451.5

$ws.Range("H132").Value = 1514.2174
$ws.Range("I132").Value = 1372.7142
$ws.Range("K132").Value = 4118.142599999999
$ws.Range("M132").Value = -1588.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1485.7
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 1307.125
$ws.Range("K113").Value = 6600
$ws.Range("L113").Value = 3921.375
$ws.Range("M113").Value = -4430
$ws.Range("N113").Value = -8261.375

$ws.Range("H122").Value = 939.8421
$ws.Range("I122").Value = 836.7143
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 7530.428699999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5080.428699999999
$ws.Range("N122").Value = -13900

$ws.Range("H137").Value = 3190.4
$ws.Range("I137").Value = 2561.7778
$ws.Range("J137").Value = 4133.3335
$ws.Range("K137").Value = 7685.3334
$ws.Range("L137").Value = 12400.0005
$ws.Range("M137").Value = -2585.3334
$ws.Range("N137").Value = -22600.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 647812.5
$ws.Range("I21").Value = 3346666.8
$ws.Range("K21").Value = 3346666.8
$ws.Range("M21").Value = -3346493.8

$ws.Range("H30").Value = 647812.5
$ws.Range("I30").Value = 3346666.8
$ws.Range("K30").Value = 3346666.8
$ws.Range("M30").Value = -3346561.8

$ws.Range("H102").Value = 4286.2
$ws.Range("I102").Value = 2216.8333
$ws.Range("J102").Value = 5665.778
$ws.Range("K102").Value = 2216.8333
$ws.Range("L102").Value = 5665.778
$ws.Range("M102").Value = -594.8332999999998
$ws.Range("N102").Value = -8909.778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("M23").Value = -1770

$ws.Range("H40").Value = 2875.5833
$ws.Range("I40").Value = 2875.5833
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2875.5833
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2739.5833
$ws.Range("N40").ClearContents()

$ws.Range("H122").Value = 4085.7585
$ws.Range("I122").Value = 3843.9583
$ws.Range("J122").Value = 5246.4
$ws.Range("K122").Value = 11531.8749
$ws.Range("L122").Value = 15739.2
$ws.Range("M122").Value = -9081.874899999999
$ws.Range("N122").Value = -20639.2

$ws.Range("H127").Value = 128559
$ws.Range("J127").Value = 128559
$ws.Range("L127").Value = 128559
$ws.Range("N127").Value = -138479

$ws.Range("H132").Value = 3235.761
$ws.Range("I132").Value = 2896.125
$ws.Range("K132").Value = 8688.375
$ws.Range("M132").Value = -6158.375

$ws.Range("H136").Value = 6548.4165
$ws.Range("I136").Value = 6052.8184
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 18158.4552
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -15608.4552
$ws.Range("N136").Value = -41100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 41866.668
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H7").Value = 15000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4887

$ws.Range("H62").Value = 7030.7
$ws.Range("I62").Value = 5826.25
$ws.Range("J62").Value = 7833.6665
$ws.Range("K62").Value = 5826.25
$ws.Range("L62").Value = 7833.6665
$ws.Range("M62").Value = -5202.25
$ws.Range("N62").Value = -9081.666499999999

$ws.Range("H65").Value = 7030.7
$ws.Range("I65").Value = 5826.25
$ws.Range("J65").Value = 7833.6665
$ws.Range("K65").Value = 29131.25
$ws.Range("L65").Value = 39168.3325
$ws.Range("M65").Value = -26011.25
$ws.Range("N65").Value = -45408.3325

$ws.Range("H132").Value = 4357.5
$ws.Range("I132").Value = 4256.316
$ws.Range("K132").Value = 12768.948
$ws.Range("M132").Value = -10238.948
